# Add the new COI (Collaborators of Interest) rows scraped from the lab partner lists
# to the bottom of the master sheet (rows 308-351), continuing the existing A/B/C
# (Last Name / First Initial / Institution) layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 308-313 continue the last styled block already present on the sheet (rows 263-307,
# font style index used there + 15pt row height). Copy that formatting down first so the
# new rows pick up the same cell style, then fix the row height to match, before filling in values.
$ws.Range("A307:C307").Copy($ws.Range("A308:C313"))
$ws.Rows("308:313").RowHeight = 15

$ws.Range("A308").Value = 'Arpaci-Dusseau'
$ws.Range("B308").Value = 'R.'
$ws.Range("C308").Value = 'U. Wisconsin-Madison'

$ws.Range("A309").Value = 'Bangalore'
$ws.Range("B309").Value = 'P.'
$ws.Range("C309").Value = 'U. of Alabama at Birmingham'

$ws.Range("A310").Value = 'Barrett'
$ws.Range("B310").Value = 'B.'
$ws.Range("C310").Value = 'Amazon'

$ws.Range("A311").Value = 'Barrett'
$ws.Range("B311").Value = 'R.'
$ws.Range("C311").Value = 'Sandia'

$ws.Range("A312").Value = 'Brandt'
$ws.Range("B312").Value = 'J.'
$ws.Range("C312").Value = 'Sandia'

$ws.Range("A313").Value = 'Brightwell'
$ws.Range("B313").Value = 'R.'
$ws.Range("C313").Value = 'Sandia'

$ws.Range("A314").Value = 'Carns'
$ws.Range("B314").Value = 'P.'
$ws.Range("C314").Value = 'ANL'

$ws.Range("A315").Value = 'Danielson'
$ws.Range("B315").Value = 'G.'
$ws.Range("C315").Value = 'Sandia'

$ws.Range("A316").Value = 'DeBardeleben'
$ws.Range("B316").Value = 'N.'
$ws.Range("C316").Value = 'LANL'

$ws.Range("A317").Value = 'Fabian'
$ws.Range("B317").Value = 'N.'
$ws.Range("C317").Value = 'Sandia'

$ws.Range("A318").Value = 'Ferreira'
$ws.Range("B318").Value = 'K.'
$ws.Range("C318").Value = 'Sandia'

$ws.Range("A319").Value = 'Gemmill'
$ws.Range("B319").Value = 'J.'
$ws.Range("C319").Value = 'Clemson U.'

$ws.Range("A320").Value = 'Gentile'
$ws.Range("B320").Value = 'A.'
$ws.Range("C320").Value = 'Sandia'

$ws.Range("A321").Value = 'Harms'
$ws.Range("B321").Value = 'K.'
$ws.Range("C321").Value = 'ANL'

$ws.Range("A322").Value = 'Harris'
$ws.Range("B322").Value = 'J.'
$ws.Range("C322").Value = 'Clemson U.'

$ws.Range("A323").Value = 'Hemmert'
$ws.Range("B323").Value = 'S.'
$ws.Range("C323").Value = 'Sandia'

$ws.Range("A324").Value = 'Kelly'
$ws.Range("B324").Value = 'S.'
$ws.Range("C324").Value = 'Sandia'

$ws.Range("A325").Value = 'Kimpe'
$ws.Range("B325").Value = 'D.'
$ws.Range("C325").Value = 'ANL'

$ws.Range("A326").Value = 'Klundt'
$ws.Range("B326").Value = 'R.'
$ws.Range("C326").Value = 'Sandia'

$ws.Range("A327").Value = 'Kroeger'
$ws.Range("B327").Value = 'T.'
$ws.Range("C327").Value = 'Sandia'

$ws.Range("A328").Value = 'Laros III'
$ws.Range("B328").Value = 'J.'
$ws.Range("C328").Value = ' Sandia'

$ws.Range("A329").Value = 'Leung'
$ws.Range("B329").Value = 'V.'
$ws.Range("C329").Value = 'Sandia'

$ws.Range("A330").Value = 'Levenhagen'
$ws.Range("B330").Value = 'M.'
$ws.Range("C330").Value = 'Sandia'

$ws.Range("A331").Value = 'Lofstead'
$ws.Range("B331").Value = 'G.'
$ws.Range("C331").Value = 'Sandia'

$ws.Range("A332").Value = 'Long'
$ws.Range("B332").Value = 'D.'
$ws.Range("C332").Value = 'UCSC'

$ws.Range("A333").Value = 'Maltzahn'
$ws.Range("B333").Value = 'C.'
$ws.Range("C333").Value = 'UCSC'

$ws.Range("A334").Value = 'Miller'
$ws.Range("B334").Value = 'E.'
$ws.Range("C334").Value = 'UCSC'

$ws.Range("A335").Value = 'Moreland'
$ws.Range("B335").Value = 'K.'
$ws.Range("C335").Value = 'Sandia'

$ws.Range("A336").Value = 'Oldfield'
$ws.Range("B336").Value = 'R.'
$ws.Range("C336").Value = 'Sandia'

$ws.Range("A337").Value = 'Pedretti'
$ws.Range("B337").Value = 'K.'
$ws.Range("C337").Value = 'Sandia'

$ws.Range("A338").Value = 'Resnick'
$ws.Range("B338").Value = 'D. R.'
$ws.Range("C338").Value = 'Sandia'

$ws.Range("A339").Value = 'Ricci'
$ws.Range("B339").Value = 'R.'
$ws.Range("C339").Value = 'U. of Utah'

$ws.Range("A340").Value = 'Rodrigues'
$ws.Range("B340").Value = 'A.'
$ws.Range("C340").Value = 'Sandia'

$ws.Range("A341").Value = 'Ross'
$ws.Range("B341").Value = 'R.'
$ws.Range("C341").Value = 'ANL'

$ws.Range("A342").Value = 'Shen'
$ws.Range("B342").Value = 'H.'
$ws.Range("C342").Value = 'Clemson U.'

$ws.Range("A343").Value = 'Skjellum'
$ws.Range("B343").Value = 'A.'
$ws.Range("C343").Value = 'Auburn U.'

$ws.Range("A344").Value = 'Sun'
$ws.Range("B344").Value = 'W.'
$ws.Range("C344").Value = 'Google'

$ws.Range("A345").Value = 'Sun'
$ws.Range("B345").Value = 'Z.'
$ws.Range("C345").Value = 'Data Direct Networks'

$ws.Range("A346").Value = 'Thompson'
$ws.Range("B346").Value = 'D.'
$ws.Range("C346").Value = 'Sandia'

$ws.Range("A347").Value = 'Tucker'
$ws.Range("B347").Value = 'T.'
$ws.Range("C347").Value = 'Sandia'

$ws.Range("A348").Value = 'Van Dyke'
$ws.Range("B348").Value = 'J.'
$ws.Range("C348").Value = ' Sandia'

$ws.Range("A349").Value = 'Vaughan'
$ws.Range("B349").Value = 'C.'
$ws.Range("C349").Value = 'Sandia'

$ws.Range("A350").Value = 'Ward'
$ws.Range("B350").Value = 'L.'
$ws.Range("C350").Value = 'Sandia'

$ws.Range("A351").Value = 'Wheeler'
$ws.Range("B351").Value = 'K.'
$ws.Range("C351").Value = 'Micron Technologies'

# Match the author's final cursor position/selection recorded in the workbook.
$ws.Range("H278").Select()
